$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename shared strings "Reflection 4" / "Reflection 5" to "Homework Reflection 4" / "Homework Reflection 5"
# and push back homework reflection 4 by one week:
#   F29 ("Reflection 4")  -> cleared
#   F31 (blank)            -> "Homework Reflection 4"
#   F32 ("Reflection 5")   -> "Homework Reflection 5"

$ws.Range("F29").Value = ""
$ws.Range("F31").Value = "Homework Reflection 4"
$ws.Range("F32").Value = "Homework Reflection 5"

# Column F must widen to fit the new, longer text
$ws.Columns("F").ColumnWidth = 19.42

# Update the selection to match final state
$ws.Range("F33").Select()
